# Commit: "Restored from revision of admin on 01/29/2021 07:00:00 AM.TEST Author: admin. Type: SAVE."
# Functional change: cell C10 on the "Rules" sheet changes from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
